$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table currently lives at A3:D13; move it to E13:H23
# (shift +4 columns, +10 rows) while keeping the header's bold style.
$src = $ws.Range("A3:D13")
$dst = $ws.Range("E13:H23")

$dst.Value2 = $src.Value2
$src.Clear() | Out-Null

# re-apply the bold header formatting to the new header row
$ws.Range("E13:H13").Font.Bold = $true

# reflect the new first/last data cell in the sheet's selection
$ws.Range("D7").Select() | Out-Null
